$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price snapshot (2021-11-09, serial 44509) takes the place of
# what used to be rows 138-139 (2021-07-06, serial 44383).
$ws.Range("D138").Value = 44509
$ws.Range("D139").Value = 44509

# The 2021-07-06 (44383) records move down into what used to be the
# 2021-08-25 (44433) slots (rows 140-141).
$ws.Range("D140").Value = 44383
$ws.Range("D141").Value = 44383

# The previous 2021-08-25 (44433) records are appended as new rows 142-143.
$ws.Range("A142").Value = 11
$ws.Range("B142").Value = "Vega Monumental Concepción"
$ws.Range("C142").Value = "Bíobío"
$ws.Range("D142").Value = 44433
$ws.Range("D142").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E142").Value = 8
$ws.Range("F142").Value = 100112040
$ws.Range("G142").Value = "Cilantro"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 200
$ws.Range("K142").Value = 600
$ws.Range("L142").Value = 700
$ws.Range("M142").Value = 650
$ws.Range("N142").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O142").Value = "Región de Ñuble"
$ws.Range("P142").Value = 650
$ws.Range("Q142").Value = 1
$ws.Range("R142").Value = "Hortaliza"

$ws.Range("A143").Value = 11
$ws.Range("B143").Value = "Vega Monumental Concepción"
$ws.Range("C143").Value = "Bíobío"
$ws.Range("D143").Value = 44433
$ws.Range("D143").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E143").Value = 8
$ws.Range("F143").Value = 100112040
$ws.Range("G143").Value = "Cilantro"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Segunda"
$ws.Range("J143").Value = 100
$ws.Range("K143").Value = 500
$ws.Range("L143").Value = 500
$ws.Range("M143").Value = 500
$ws.Range("N143").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O143").Value = "Región de Ñuble"
$ws.Range("P143").Value = 500
$ws.Range("Q143").Value = 1
$ws.Range("R143").Value = "Hortaliza"
